$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update res_bus vm_pu results for the 380 kV case (slack bus voltage 1.05 -> 1.02 p.u.,
# with corresponding recalculated per-unit voltages for all buses/rows 2-25).
$values = @{
    "B2" = 1.02
    "C2" = 1.040393660519047
    "D2" = 1.052731842878131
    "E2" = 1.048992874139455
    "F2" = 1.060435114373241
    "I2" = 1.038162686421549
    "J2" = 1.045480507182045
    "K2" = 1.055479762953657
    "L2" = 1.051751172556009
    "M2" = 1.063161902151157
    "N2" = 1.018898674364074
    "B3" = 1.02
    "C3" = 1.041867824937321
    "D3" = 1.053544703387799
    "E3" = 1.050238236439371
    "F3" = 1.061612165095181
    "I3" = 1.038332649511734
    "J3" = 1.046597652730279
    "K3" = 1.056105527734218
    "L3" = 1.052807572192282
    "M3" = 1.064152462235481
    "N3" = 1.019283512923265
    "B4" = 1.02
    "C4" = 1.042821091374123
    "D4" = 1.054069511238092
    "E4" = 1.051043632110217
    "F4" = 1.062373030382255
    "I4" = 1.038440833770517
    "J4" = 1.047319514645559
    "K4" = 1.056508622750124
    "L4" = 1.053490138980301
    "M4" = 1.064792080327734
    "N4" = 1.01953182650147
    "B5" = 1.02
    "C5" = 1.043221702675756
    "D5" = 1.05428986177078
    "E5" = 1.051382119395412
    "F5" = 1.062692717629232
    "I5" = 1.038485885348872
    "J5" = 1.047622748418013
    "K5" = 1.056677649881117
    "L5" = 1.053776855249788
    "M5" = 1.065060657300485
    "N5" = 1.019636050493136
    "B6" = 1.02
    "C6" = 1.043288958919428
    "D6" = 1.054326843274073
    "E6" = 1.051438947077474
    "F6" = 1.062746384016491
    "I6" = 1.038493424551629
    "J6" = 1.047673648917903
    "K6" = 1.056706004855091
    "L6" = 1.053824982498157
    "M6" = 1.065105733951474
    "N6" = 1.019653540399358
    "B7" = 1.02
    "C7" = 1.04282644491249
    "D7" = 1.05407245666686
    "E7" = 1.051048155386552
    "F7" = 1.062377302764921
    "I7" = 1.03844143743767
    "J7" = 1.04732356739646
    "K7" = 1.056510883003236
    "L7" = 1.053493971015741
    "M7" = 1.064795670314767
    "N7" = 1.019533219802413
    "B8" = 1.02
    "C8" = 1.040891992368139
    "D8" = 1.053006795040693
    "E8" = 1.049413842797279
    "F8" = 1.060833064241002
    "I8" = 1.038220497676018
    "J8" = 1.045858262362809
    "K8" = 1.055691619781082
    "L8" = 1.05210839562438
    "M8" = 1.063496945612151
    "N8" = 1.019028878575625
    "B9" = 1.02
    "C9" = 1.037478263005937
    "D9" = 1.051119988241249
    "E9" = 1.046530467826436
    "F9" = 1.058105936928422
    "I9" = 1.037817429331266
    "J9" = 1.043268330267437
    "K9" = 1.054234018927504
    "L9" = 1.049659072231642
    "M9" = 1.061198049382963
    "N9" = 1.018134729341642
    "B10" = 1.02
    "C10" = 1.035198739282798
    "D10" = 1.049856030061594
    "E10" = 1.044605633152152
    "F10" = 1.056283658148425
    "I10" = 1.037539460253379
    "J10" = 1.041536168864759
    "K10" = 1.053252842331894
    "L10" = 1.048020763976659
    "M10" = 1.059658317285326
    "N10" = 1.017534899242527
    "B11" = 1.02
    "C11" = 1.034210721676216
    "D11" = 1.049307265565162
    "E11" = 1.043771494328792
    "F11" = 1.055493557342103
    "I11" = 1.037416896905845
    "J11" = 1.040784755059387
    "K11" = 1.052825727908975
    "L11" = 1.047310024754021
    "M11" = 1.058989868143454
    "N11" = 1.017274264108796
    "B12" = 1.02
    "C12" = 1.033843575427478
    "D12" = 1.049103209067178
    "E12" = 1.043461552472982
    "F12" = 1.055199918671496
    "I12" = 1.037371040338728
    "J12" = 1.040505434914882
    "K12" = 1.052666737919929
    "L12" = 1.047045818991387
    "M12" = 1.058741312503428
    "N12" = 1.017177315161725
    "B13" = 1.02
    "C13" = 1.03392233655819
    "D13" = 1.049146989896405
    "E13" = 1.043528040902242
    "F13" = 1.055262912462175
    "I13" = 1.037380891717085
    "J13" = 1.040565359686773
    "K13" = 1.052700857250645
    "L13" = 1.047102501393055
    "M13" = 1.058794640524166
    "N13" = 1.01719811728466
    "B14" = 1.02
    "C14" = 1.034180376404526
    "D14" = 1.049290402705554
    "E14" = 1.043745876617704
    "F14" = 1.0554692883633
    "I14" = 1.037413113144158
    "J14" = 1.040761670713552
    "K14" = 1.05281259270451
    "L14" = 1.047288189640682
    "M14" = 1.058969327858125
    "N14" = 1.017266253091466
    "B15" = 1.02
    "C15" = 1.034339342851016
    "D15" = 1.049378734757566
    "E15" = 1.043880078269377
    "F15" = 1.055596422053865
    "I15" = 1.037432921944389
    "J15" = 1.040882596233072
    "K15" = 1.052881391421487
    "L15" = 1.047402570935182
    "M15" = 1.059076923428881
    "N15" = 1.017308215551026
    "B16" = 1.02
    "C16" = 1.035264289752142
    "D16" = 1.049892418840206
    "E16" = 1.044660977642916
    "F16" = 1.056336072284097
    "I16" = 1.037547547972775
    "J16" = 1.041586008314333
    "K16" = 1.053281140809892
    "L16" = 1.048067904795969
    "M16" = 1.05970264322635
    "N16" = 1.017552177546053
    "B17" = 1.02
    "C17" = 1.03584422022522
    "D17" = 1.050214246934122
    "E17" = 1.045150631924554
    "F17" = 1.056799753882845
    "I17" = 1.037618860312269
    "J17" = 1.042026868224827
    "K17" = 1.053531287423098
    "L17" = 1.048484889248932
    "M17" = 1.06009467396498
    "N17" = 1.017704965160524
    "B18" = 1.02
    "C18" = 1.036182390857445
    "D18" = 1.050401822815365
    "E18" = 1.045436174381919
    "F18" = 1.057070111208398
    "I18" = 1.037660243267602
    "J18" = 1.042283881867126
    "K18" = 1.053676975893755
    "L18" = 1.048727980124275
    "M18" = 1.060323171680263
    "N18" = 1.017793996400547
    "B19" = 1.02
    "C19" = 1.036297682739239
    "D19" = 1.050465757464027
    "E19" = 1.045533526158401
    "F19" = 1.057162279160895
    "I19" = 1.037674317784269
    "J19" = 1.042371494550038
    "K19" = 1.053726614992416
    "L19" = 1.048810846056492
    "M19" = 1.060401055235464
    "N19" = 1.017824339009463
    "B20" = 1.02
    "C20" = 1.035782008844023
    "D20" = 1.050179732418937
    "E20" = 1.045098103368823
    "F20" = 1.056750015636144
    "I20" = 1.0376112311381
    "J20" = 1.0419795818671
    "K20" = 1.053504471598666
    "L20" = 1.048440164146881
    "M20" = 1.06005263009264
    "N20" = 1.017688581525095
    "B21" = 1.02
    "C21" = 1.034104394367651
    "D21" = 1.049248177316566
    "E21" = 1.043681732379001
    "F21" = 1.055408520250604
    "I21" = 1.037403633879914
    "J21" = 1.040703867888979
    "K21" = 1.052779698813627
    "L21" = 1.047233514787775
    "M21" = 1.058917894102096
    "N21" = 1.017246192601306
    "B22" = 1.02
    "C22" = 1.033048725660299
    "D22" = 1.048661192784879
    "E22" = 1.042790590758119
    "F22" = 1.054564141096988
    "I22" = 1.037271193486909
    "J22" = 1.039900549225213
    "K22" = 1.052322033666761
    "L22" = 1.046473655026667
    "M22" = 1.058202912294017
    "N22" = 1.016967249098488
    "B23" = 1.02
    "C23" = 1.033608441924523
    "D23" = 1.048972486032536
    "E23" = 1.043263061400003
    "F23" = 1.055011851575952
    "I23" = 1.037341584379447
    "J23" = 1.040326521529508
    "K23" = 1.052564837977508
    "L23" = 1.046876585373671
    "M23" = 1.058582083614892
    "N23" = 1.01711519829346
    "B24" = 1.02
    "C24" = 1.035810119777039
    "D24" = 1.050195328478798
    "E24" = 1.045121838962098
    "F24" = 1.056772490519255
    "I24" = 1.03761467908973
    "J24" = 1.04200094894643
    "K24" = 1.053516589188829
    "L24" = 1.048460373894581
    "M24" = 1.060071628425759
    "N24" = 1.01769598485427
    "B25" = 1.02
    "C25" = 1.038361420407233
    "D25" = 1.051608841952298
    "E25" = 1.047276328250823
    "F25" = 1.058811691210433
    "I25" = 1.037923262542851
    "J25" = 1.043938848313641
    "K25" = 1.054612503521149
    "L25" = 1.050293222161848
    "M25" = 1.061793615065696
    "N25" = 1.0183665402894
}

foreach ($key in $values.Keys) {
    $ws.Range($key).Value = $values[$key]
}
